# Implementacion temporal de GLOSAS CARGUE
$wb = $excel.ActiveWorkbook
$glosas = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Create the new (hidden, helper) DETAILS sheet right after GLOSAS
# ---------------------------------------------------------------------------
$details = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $glosas)
$details.Name = "DETAILS"

# ---------------------------------------------------------------------------
# 2. Populate DETAILS content
# ---------------------------------------------------------------------------
$details.Range("A1").Value = "VERSION"
$details.Range("A2").Value = "TIPO DE CARGA"
$details.Range("A3").Value = "HUB"

# "1.0" must land as literal text (not get coerced to the number 1), so we
# stage the cells as Text before typing, then restore the date display
# format the original column already carried.
$details.Range("B1:B3").NumberFormat = "@"
$details.Range("B1").Value = "1.0"
$details.Range("B2").Value = "MANUAL"
$details.Range("B3").Value = "Fact.Squid"
$details.Range("B1:B3").NumberFormat = "m/d/yy"

$details.Range("A1:A3").Font.Bold = $true
$details.Range("A1:A3").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$details.Range("A1:A3").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

$details.Range("A5").Value = "TIPOS"
$details.Range("A5").Font.Bold = $true

$details.Range("B5").Value = "0-GLOSA_INICIAL"
$details.Range("B6").Value = "1-GLOSA_RATIFICADA"
$details.Range("B7").Value = "2-GLOSA_INICIAL_SALDO_CERO"
$details.Range("B8").Value = "3-GLOSA_RATIFICADA_SALDO_CERO"
$details.Range("B9").Value = "4-DEVOLUCION"
$details.Range("B10").Value = "5-FACTURA_DEVUELTA(GLOSA 100%)"
$details.Range("B11").Value = "6-FACTURA_DEVUELTA(NO_ACEPTACION)"
$details.Range("B5:B11").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$details.Columns.Item(1).AutoFit()
$details.Columns.Item(2).AutoFit()

$details.Range("B18:B19").Select()

# ---------------------------------------------------------------------------
# 3. Rebuild GLOSAS: push the original header row down to row 6 and add the
#    new "document" header block in rows 1-4.
# ---------------------------------------------------------------------------
$glosas.Range("1:5").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

$glosas.Range("A1").Value = "DOCUMENTO GESTOR"
$glosas.Range("A2").Value = "FECHA DOCUMENTO"
$glosas.Range("A3").Value = "TIPO"
$glosas.Range("A4").Value = "NIT"

$glosas.Range("A1:A4").Font.Bold = $true
$glosas.Range("A1:A4").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$glosas.Range("A1:A4").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Column B already carries a date style by default (inherited from the
# original workbook). Strip that back to plain "Normal" on the fields that
# are NOT dates, and keep/ (re)apply the date format on the date field (B2)
# and its merged continuation (C2:F2).
$glosas.Range("B1").Style = "Normal"
$glosas.Range("B3").Style = "Normal"
$glosas.Range("B4").Style = "Normal"
$glosas.Range("B2:F2").NumberFormat = "m/d/yy"

$glosas.Range("B1:F1").Merge()
$glosas.Range("B2:F2").Merge()
$glosas.Range("B3:F3").Merge()
$glosas.Range("B4:F4").Merge()

$glosas.Range("B1:F4").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$glosas.Range("B1").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$glosas.Range("B1").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$glosas.Range("B2").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$glosas.Range("B2").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Header row (now row 6) gets bold + border like the other labels
$glosas.Range("A6:F6").Font.Bold = $true
$glosas.Range("A6:F6").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$glosas.Range("A6:F6").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# ---------------------------------------------------------------------------
# 4. Data validation dropdown on B3 (TIPO) sourced from DETAILS!B5:B11
# ---------------------------------------------------------------------------
$glosas.Range("B3").Validation.Add(3, 1, 1, "=DETAILS!`$B`$5:`$B`$11")
$glosas.Range("B3").Validation.IgnoreBlank = $true
$glosas.Range("B3").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------------
# 5. Column widths / layout touch-ups
# ---------------------------------------------------------------------------
$glosas.Columns.Item(1).ColumnWidth = 24.6
$glosas.Columns.Item(3).ColumnWidth = 15.7
$glosas.Columns.Item(4).ColumnWidth = 18.9
$glosas.Columns.Item(5).ColumnWidth = 23.7
$glosas.Columns.Item(6).AutoFit()

$glosas.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# ---------------------------------------------------------------------------
# 6. Selection / visibility
# ---------------------------------------------------------------------------
$glosas.Activate()
$glosas.Range("B1:F1").Select()

$details.Visible = [Microsoft.Office.Interop.Excel.XlSheetVisibility]::xlSheetHidden

Write-Host "done"
